$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Cells.Item(2, 6)
Write-Host "NumberFormat: $($c.NumberFormat)"
Write-Host "Font.Bold: $($c.Font.Bold)"
Write-Host "Font.Size: $($c.Font.Size)"
Write-Host "Interior.Color: $($c.Interior.Color)"
Write-Host "Borders(xlEdgeBottom).LineStyle: $($c.Borders.Item(9).LineStyle)"
